# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the
# 5860fed6-83c1-4577-8548-15d0a1a18605 record on both the
# "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-29 02:18:29"
$wsZhCn.Range("G4").Value = "2016-01-29 02:19:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-29 02:18:42"
$wsDeDe.Range("G4").Value = "2016-01-29 02:19:39"
